$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "41.213.32"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.181.01"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.81%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "251.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "

# Row 6
$ws.Range("E6").Value = "  -3.09%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "66.66"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -6.90%  "

# Row 8
$ws.Range("E8").Value = "  +0.10%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.590"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "59.09"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "36.61"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -9.77%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0937"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.90%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.97%  "

# Row 14
$ws.Range("E14").Value = "  -4.84%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.506.45"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.31"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -4.28%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.846"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.52%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.173.42"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.05%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "41.143.20"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.19%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0947"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.78%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "71.65"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.06"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.61%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "230.42"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.73%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.82"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.18%  "

# Row 26
$ws.Range("E26").Value = "  +0.29%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.39"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.43%  "

# Row 28
$ws.Range("E28").Value = "  -4.65%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "168.03"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.71%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.65%  "

# Row 31
$ws.Range("E31").Value = "  -2.75%  "

# Row 32
$ws.Range("E32").Value = "  -1.96%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.71"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.64%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0749"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.09%  "

# Row 35
$ws.Range("E35").Value = "  -2.05%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.51"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.55%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.97"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "24.60"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -5.60%  "

# Row 39
$ws.Range("E39").Value = "  +0.23%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.56"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +14.60%  "

# Row 41
$ws.Range("E41").Value = "  -3.48%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.53"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -7.23%  "

# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "60.96"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -8.11%  "

# Row 44
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "11.34"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -7.58%  "

# Row 45
$ws.Range("E45").Value = "  -1.87%  "

# Row 46
$ws.Range("E46").Value = "  -6.71%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0993"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.25%  "

# Row 48
$ws.Range("E48").Value = "  -0.20%  "

# Row 49
$ws.Range("E49").Value = "  -2.49%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.24"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -9.81%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.15"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.37%  "
